$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.180.47'
$ws.Range('E2').Value = '  -5.70%  '
$ws.Range('D3').Value = '3.705.26'
$ws.Range('E3').Value = '  -5.39%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = "'586.47"
$ws.Range('E5').Value = '  -2.60%  '
$ws.Range('D6').Value = "'181.40"
$ws.Range('E6').Value = '  +7.41%  '
$ws.Range('D7').Value = '3.702.13'
$ws.Range('E7').Value = '  -5.29%  '
$ws.Range('D8').Value = "'0.629"
$ws.Range('E8').Value = '  -6.67%  '
$ws.Range('E9').Value = '  -0.14%  '
$ws.Range('D10').Value = "'0.714"
$ws.Range('E10').Value = '  -7.60%  '
$ws.Range('D11').Value = "'0.163"
$ws.Range('E11').Value = '  -10.51%  '
$ws.Range('D12').Value = "'54.54"
$ws.Range('E12').Value = '  -0.37%  '
$ws.Range('E13').Value = '  -10.68%  '
$ws.Range('E14').Value = '  -8.69%  '
$ws.Range('D15').Value = '4.291.09'
$ws.Range('E15').Value = '  -5.16%  '
$ws.Range('D16').Value = '3.702.86'
$ws.Range('E16').Value = '  -5.16%  '
$ws.Range('D17').Value = "'19.42"
$ws.Range('E17').Value = '  -8.60%  '
$ws.Range('E18').Value = '  -2.96%  '
$ws.Range('D19').Value = "'12.85"
$ws.Range('E19').Value = '  -8.77%  '
$ws.Range('E20').Value = '  -8.17%  '
$ws.Range('D21').Value = '67.760.26'
$ws.Range('E21').Value = '  -5.87%  '
$ws.Range('D22').Value = "'408.25"
$ws.Range('E22').Value = '  -7.18%  '
$ws.Range('D23').Value = "'4.50"
$ws.Range('E23').Value = '  -5.91%  '
$ws.Range('D24').Value = "'88.48"
$ws.Range('E24').Value = '  -6.69%  '
$ws.Range('D25').Value = "'3.04"
$ws.Range('E25').Value = '  -8.13%  '
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Value = "'12.77"
$ws.Range('E26').Value = '  -8.78%  '
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D27').Value = "'10.99"
$ws.Range('E27').Value = '  -0.83%  '
$ws.Range('E28').Value = '  -8.26%  '
$ws.Range('E29').Value = '  +2.12%  '
$ws.Range('D30').Value = "'9.52"
$ws.Range('E30').Value = '  -7.61%  '
$ws.Range('E31').Value = '  -8.11%  '
$ws.Range('E32').Value = '  -4.86%  '
$ws.Range('D33').Value = "'12.48"
$ws.Range('E33').Value = '  -9.18%  '
$ws.Range('E34').Value = '  -8.12%  '
$ws.Range('E35').Value = '  -4.71%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = "'43.33"
$ws.Range('E36').Value = '  -15.23%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').Value = "'601.74"
$ws.Range('E37').Value = '  -1.63%  '
$ws.Range('D38').Value = '0.0₃0896'
$ws.Range('E38').Value = '  -10.05%  '
$ws.Range('E39').Value = '  -0.05%  '
$ws.Range('D40').Value = "'0.398"
$ws.Range('E40').Value = '  -6.18%  '
$ws.Range('E41').Value = '  +0.18%  '
$ws.Range('E42').Value = '  -5.05%  '
$ws.Range('D43').Value = "'2.79"
$ws.Range('E43').Value = '  +5.30%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = "'2.99"
$ws.Range('E44').Value = '  -10.52%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = "'2.98"
$ws.Range('E45').Value = '  -7.93%  '
$ws.Range('D46').Value = "'0.0434"
$ws.Range('E46').Value = '  -8.43%  '
$ws.Range('D47').Value = "'9.25"
$ws.Range('E47').Value = '  -10.89%  '
$ws.Range('D48').Value = '2.804.55'
$ws.Range('E48').Value = '  -2.88%  '
$ws.Range('E49').Value = '  -8.15%  '
$ws.Range('D50').Value = "'2.69"
$ws.Range('E50').Value = '  -5.72%  '
$ws.Range('D51').Value = "'3.09"
$ws.Range('E51').Value = '  -7.48%  '
